# "Refactor code for run config"
# Merge the VIE/ENG duplicate test-case blocks on Sheet1 into a single
# block with separate "Result VIE" / "Result ENG" columns, renumber the
# TestCase IDs (UP_SUPPLIER_VIE_0X -> UP_SUPPLIER_0X), drop the now
# redundant "(VIE)"/"(ENG)" language qualifiers from the Expected text,
# remove the old ENG-only rows (10:17), and shrink the print area back
# down to the real data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Header row -----------------------------------------------------
$ws.Range('F1').Value = 'Result VIE'
$ws.Range('G1').Value = 'Result ENG'
$ws.Range('A1').Copy()
$ws.Range('G1').PasteSpecial(-4122)   # xlPasteFormats, keep header style

# --- Row 2 : UP_SUPPLIER_01 (View supplier detail) ------------------
$ws.Range('A2').Value = 'UP_SUPPLIER_01'
$ws.Range('G2').Value = $ws.Range('F2').Text
$ws.Range('G2').Style = 'Normal'

# --- Row 3 : UP_SUPPLIER_02 (Update VN supplier) ---------------------
$ws.Range('A3').Value = 'UP_SUPPLIER_02'
$ws.Range('E3').Value = '1. Text must be translated correctly according to the console language ' + $nl + '2. Can update supplier information'
$ws.Range('G3').Value = $ws.Range('F3').Text
$ws.Range('G3').Style = 'Normal'

# --- Row 4 : UP_SUPPLIER_03 (Update Non-VN supplier) -----------------
$ws.Range('A4').Value = 'UP_SUPPLIER_03'
$ws.Range('E4').Value = '1. Text must be translated correctly according to the console language ' + $nl + '2. Can update supplier information'
$ws.Range('G4').Value = $ws.Range('F4').Text
$ws.Range('G4').Style = 'Normal'

# --- Row 5 : UP_SUPPLIER_04 (blank required field) -------------------
$ws.Range('A5').Value = 'UP_SUPPLIER_04'
$ws.Range('E5').Value = '1. Error should be shown and  must be translated correctly according to the console language ' + $nl + '2. Can not Update supplier'
$ws.Range('G5').Value = $ws.Range('F5').Text
$ws.Range('G5').Style = 'Normal'

# --- Row 6 : UP_SUPPLIER_05 (available supplier code) ----------------
$ws.Range('A6').Value = 'UP_SUPPLIER_05'
$ws.Range('E6').Value = '1. Error should be shown and  must be translated correctly according to the console language ' + $nl + '2. Can not Update supplier'
$ws.Range('G6').Value = $ws.Range('F6').Text
$ws.Range('G6').Style = 'Normal'

# --- Row 7 : UP_SUPPLIER_06 (invalid format supplier code) -----------
$ws.Range('A7').Value = 'UP_SUPPLIER_06'
$ws.Range('E7').Value = '1. Error should be shown and  must be translated correctly according to the console language ' + $nl + '2. Can not update supplier'
$ws.Range('G7').Value = $ws.Range('F7').Text
$ws.Range('G7').Style = 'Normal'

# --- Row 8 : UP_SUPPLIER_07 (Check Order history) ---------------------
$ws.Range('A8').Value = 'UP_SUPPLIER_07'
$ws.Range('G8').Value = $ws.Range('F8').Text
$ws.Range('G8').Style = 'Normal'

# --- Row 9 : UP_SUPPLIER_08 (Delete Supplier) -------------------------
$ws.Range('A9').Value = 'UP_SUPPLIER_08'
$ws.Range('G9').Value = $ws.Range('F9').Text
$ws.Range('G9').Style = 'Normal'

# --- Drop the old English-only duplicate rows -------------------------
$ws.Rows('10:17').Delete()

# --- Shrink the print area back to the live data range -----------------
$ws.PageSetup.PrintArea = '$A$1:$G$9'
